$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text interpretation so numeric-looking strings (e.g. "248.76")
    # are not silently parsed into floating point numbers, then restore the
    # original (default) cell style so no formatting change is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue 2  4 "248.76"
Set-TextValue 3  4 "22.71"
Set-TextValue 4  4 "5.273"
Set-TextValue 5  4 "0.05696"
Set-TextValue 7  4 "6.332"
Set-TextValue 8  4 "0.8057"
Set-TextValue 9  4 "0.8959"
Set-TextValue 10 4 "0.1427"
Set-TextValue 11 4 "0.07455"
Set-TextValue 12 4 "0.03098"
Set-TextValue 14 4 "0.09402"
Set-TextValue 15 4 "3.862"
Set-TextValue 16 4 "0.001579"
Set-TextValue 17 4 "0.04806"
Set-TextValue 18 4 "0.01828"
Set-TextValue 19 4 "0.0005808"
Set-TextValue 20 4 "0.006425"
Set-TextValue 21 4 "0.004991"
Set-TextValue 22 4 "0.0009972"
Set-TextValue 23 4 "0.0001501"
Set-TextValue 24 4 "3.695"
Set-TextValue 27 4 "0.1369"
Set-TextValue 40 4 "0.03980"
Set-TextValue 41 4 "0.006816"
Set-TextValue 42 4 "0.1068"
Set-TextValue 43 4 "0.002751"
Set-TextValue 44 4 "0.007695"
Set-TextValue 45 4 "0.00005589"
Set-TextValue 46 4 "0.00000000751"
Set-TextValue 47 4 "0.4989"
Set-TextValue 48 4 "0.2008"
Set-TextValue 49 4 "0.00002101"
Set-TextValue 50 4 "0.01011"

# Column E (Volume(1h)) updates - plain text, no numeric coercion risk
$ws.Cells.Item(19, 5).Value = "18OneONEWorstin24h"
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOIN"
